# Added posts export functionality:
#  - Extend the slug formula in column G so that, in addition to spaces,
#    colons, semicolons, commas, periods and slashes are also replaced
#    with a hyphen before lower-casing (keeps generated slugs URL-safe).
#  - Leave the cursor/selection where the author last left it while
#    making the change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("posts")
$ws.Activate()

$slugFormula = '=LOWER(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(B2, " ", "-"), ":", "-"), ";", "-"), ",", "-"), ".", "-"), "/", "-"))'

# G2 keeps its own (non-shared) formula.
$ws.Range("G2").Formula = $slugFormula

# G3:G11 form the shared-formula block anchored at G3; assigning the same
# formula text across the whole range lets Excel re-derive it as a shared
# formula (ref="G3:G11", si="0") exactly like the original block.
$ws.Range("G3:G11").Formula = '=LOWER(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(B3, " ", "-"), ":", "-"), ";", "-"), ",", "-"), ".", "-"), "/", "-"))'

# Move the view: scroll so column B is the left-most visible column and
# select F15 (matches the author's cursor position when they saved).
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("F15").Select()
